$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$newXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Title"/><w:spacing w:line="360" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">Binary Search </w:t></w:r><w:r><w:t>Tree</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
try {
  $p1.Range.InsertXML($newXml)
  Write-Host "InsertXML worked"
} catch {
  Write-Host "InsertXML failed: $_"
}
